# Updated cryptos list on Sun May  7 17:35:02 UTC 2023 with GitHub Actions
# Applies latest price/volume(1h) snapshot to the cryptos sheet.
# Cells whose new text would otherwise be auto-parsed as a number
# (single-decimal-point values like "1.005") are forced to Text
# format first so they are stored as literal strings, matching the
# original inline-string cell contents (e.g. "29.045.51",
# "1.922.87" already survive as text because of the double dot).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.045.51'
$ws.Range("E2").Value = '  +0.67%  '

$ws.Range("D3").Value = '1.922.87'
$ws.Range("E3").Value = '  +1.76%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.34%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.58'
$ws.Range("E5").Value = '  +0.85%  '

$ws.Range("E6").Value = '  +0.32%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4595'
$ws.Range("E7").Value = '  +0.57%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3819'
$ws.Range("E8").Value = '  +0.49%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07753'
$ws.Range("E9").Value = '  +0.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9796'
$ws.Range("E10").Value = '  +1.67%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.68'
$ws.Range("E11").Value = '  +2.88%  '

$ws.Range("D12").Value = '1.911.66'
$ws.Range("E12").Value = '  +1.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.690'
$ws.Range("E13").Value = '  +0.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.964'
$ws.Range("E14").Value = '  +0.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07030'
$ws.Range("E15").Value = '  -0.01%  '

$ws.Range("E16").Value = '  +0.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.22'
$ws.Range("E17").Value = '  +1.28%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009518'
$ws.Range("E18").Value = '  +0.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.72'
$ws.Range("E19").Value = '  +0.67%  '

$ws.Range("E20").Value = '  +0.28%  '

$ws.Range("D21").Value = '29.053.90'
$ws.Range("E21").Value = '  +0.92%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.338'
$ws.Range("E22").Value = '  +0.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.97'
$ws.Range("E23").Value = '  +1.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.076'
$ws.Range("E24").Value = '  +0.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.32'
$ws.Range("E25").Value = '  +1.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.09'
$ws.Range("E26").Value = '  +0.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.658'
$ws.Range("E27").Value = '  +1.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.94'
$ws.Range("E28").Value = '  +0.70%  '

$ws.Range("E29").Value = '  +1.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09340'
$ws.Range("E30").Value = '  +0.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.8577'
$ws.Range("E31").Value = '  +1.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.106'
$ws.Range("E32").Value = '  +0.78%  '

$ws.Range("E33").Value = '  +0.46%  '

$ws.Range("E34").Value = '  +0.27%  '

$ws.Range("B35").Value = 'TrustWalletToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.160'
$ws.Range("E35").Value = '  +1.46%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05682'
$ws.Range("E36").Value = '  +0.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.166'
$ws.Range("E37").Value = '  +17.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.004'
$ws.Range("E38").Value = '  +0.39%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02046'
$ws.Range("E39").Value = '  +0.87%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.482'
$ws.Range("E40").Value = '  +1.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5503'
$ws.Range("E41").Value = '  +0.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1754'
$ws.Range("E42").Value = '  +0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.364'
$ws.Range("E43").Value = '  +2.52%  '

$ws.Range("E44").Value = '  +6.39%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000002772'
$ws.Range("E45").Value = '  -7.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5181'
$ws.Range("E46").Value = '  +0.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.26'
$ws.Range("E47").Value = '  -0.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06911'
$ws.Range("E48").Value = '  +1.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '110.27'
$ws.Range("E49").Value = '  -1.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.765'
$ws.Range("E50").Value = '  -0.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.004'
$ws.Range("E51").Value = '  +0.41%  '

